$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 215, pushing existing rows 215-219 down to 216-220
$ws.Rows.Item(215).Insert()

# Populate the new row 215 with a new weekly record (same dimension/category values,
# new date and new price observations)
$ws.Cells.Item(215, 1).Value = 10
$ws.Cells.Item(215, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(215, 3).Value = "La Araucanía"
$ws.Cells.Item(215, 4).Value = 44656
$ws.Cells.Item(215, 5).Value = 9
$ws.Cells.Item(215, 6).Value = "Fruta"
$ws.Cells.Item(215, 7).Value = 100102
$ws.Cells.Item(215, 8).Value = "Cítricos"
$ws.Cells.Item(215, 9).Value = 100102006
$ws.Cells.Item(215, 10).Value = "Pomelo"
$ws.Cells.Item(215, 11).Value = "Start Ruby"
$ws.Cells.Item(215, 12).Value = "Especial"
$ws.Cells.Item(215, 13).Value = 20
$ws.Cells.Item(215, 14).Value = 18000
$ws.Cells.Item(215, 15).Value = 18000
$ws.Cells.Item(215, 16).Value = 18000
$ws.Cells.Item(215, 17).Value = "$/bandeja 15 kilos granel"
$ws.Cells.Item(215, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(215, 19).Value = 1200
$ws.Cells.Item(215, 20).Value = 15
